# Updated cryptos list on Fri Feb  9 09:09:55 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    # Force the cell to hold the exact literal text, even when the text
    # looks like a number (e.g. "323.30" or "1.00"), by temporarily
    # switching the cell to Text format before assigning the value, then
    # restoring the default ("Normal") style so no stray formatting is
    # left behind on the cell.
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "46.620.26"
Set-TextValue "E2" "  +4.31%  "

Set-TextValue "D3" "2.474.31"
Set-TextValue "E3" "  +2.03%  "

Set-TextValue "E4" "  +0.00%  "

Set-TextValue "D5" "323.30"
Set-TextValue "E5" "  +2.48%  "

Set-TextValue "D6" "106.12"
Set-TextValue "E6" "  +3.92%  "

Set-TextValue "D7" "0.521"
Set-TextValue "E7" "  +1.53%  "

Set-TextValue "D8" "0.999"
Set-TextValue "E8" "  -0.08%  "

Set-TextValue "E9" "  +2.94%  "

Set-TextValue "D10" "36.29"
Set-TextValue "E10" "  +1.95%  "

Set-TextValue "E12" "  +0.53%  "

Set-TextValue "D13" "18.43"
Set-TextValue "E13" "  -2.75%  "

Set-TextValue "D14" "7.11"
Set-TextValue "E14" "  +2.12%  "

Set-TextValue "D15" "2.867.91"
Set-TextValue "E15" "  +2.34%  "

Set-TextValue "D16" "2.428.09"
Set-TextValue "E16" "  -0.49%  "

Set-TextValue "D17" "0.847"
Set-TextValue "E17" "  +1.46%  "

Set-TextValue "D18" "46.510.93"
Set-TextValue "E18" "  +4.36%  "

Set-TextValue "D19" "12.73"
Set-TextValue "E19" "  +2.63%  "

Set-TextValue "D20" "6.48"
Set-TextValue "E20" "  +1.15%  "

Set-TextValue "D21" "0.0₃0940"
Set-TextValue "E21" "  +2.00%  "

Set-TextValue "D22" "70.68"

Set-TextValue "D23" "2.38"
Set-TextValue "E23" "  +3.67%  "

Set-TextValue "D24" "248.94"
Set-TextValue "E24" "  +2.62%  "

Set-TextValue "D25" "2.55"
Set-TextValue "E25" "  +2.82%  "

Set-TextValue "D26" "26.23"
Set-TextValue "E26" "  +3.77%  "

Set-TextValue "E27" "  +0.00%  "

Set-TextValue "D28" "2.19"
Set-TextValue "E28" "  -3.87%  "

Set-TextValue "E29" "  +2.96%  "

Set-TextValue "D30" "34.85"
Set-TextValue "E30" "  +3.52%  "

Set-TextValue "D31" "49.68"
Set-TextValue "E31" "  +2.60%  "

Set-TextValue "E32" "  +3.21%  "

Set-TextValue "D33" "19.87"
Set-TextValue "E33" "  +2.24%  "

Set-TextValue "E34" "  +3.16%  "

# Rows 35 and 36 swap contents (Hedera and FirstDigitalUSD trade places).
Set-TextValue "B35" "Hedera"
Set-TextValue "C35" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D35" "0.0770"
Set-TextValue "E35" "  -1.10%  "

Set-TextValue "B36" "FirstDigitalUSD"
Set-TextValue "C36" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D36" "1.00"
Set-TextValue "E36" "  +0.07%  "

Set-TextValue "D37" "4.62"
Set-TextValue "E37" "  +2.92%  "

Set-TextValue "D38" "1.91"
Set-TextValue "E38" "  +1.06%  "

Set-TextValue "E39" "  +2.96%  "

Set-TextValue "D40" "124.29"
Set-TextValue "E40" "  +3.74%  "

Set-TextValue "E41" "  +2.23%  "

Set-TextValue "E42" "  +1.56%  "

Set-TextValue "D43" "21.08"
Set-TextValue "E43" "  -0.19%  "

Set-TextValue "D44" "0.0294"
Set-TextValue "E44" "  +1.27%  "

Set-TextValue "D45" "1.983.93"
Set-TextValue "E45" "  +2.08%  "

Set-TextValue "D47" "2.10"
Set-TextValue "E47" "  -2.86%  "

Set-TextValue "E48" "  +10.24%  "

Set-TextValue "D49" "9.10"
Set-TextValue "E49" "  -3.88%  "

Set-TextValue "D50" "5.17"
Set-TextValue "E50" "  +11.05%  "

Set-TextValue "D51" "79.24"
Set-TextValue "E51" "  +4.74%  "
